$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1.25
$ws.Range("D3").Value = 1.4
$ws.Range("B4").Value = 1.49
$ws.Range("E4").Value = 1.21
$ws.Range("C5").Value = 1.37
$ws.Range("D5").Value = 1.35
$ws.Range("E6").Value = 1.33
